# Add a new scenario row (row 69) to the "Example Scenarios" sheet:
#   A69 -> short description of the scenario
#   B69 -> the multi-line "statements:" rule body
#
# Mirrors the author's upload that appended a "check if the last login
# occured in the within 30 minutes" example to the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newline = [char]10

$aText = "check if the last login occured in the within 30 minutes"

$bText = "statements:" + $newline +
    "    - context: ""lastLogin := user['urn:ietf:params:scim:schemas:extension:ibm:2.0:User'].lastLogin""" + $newline +
    "    - context: ""currentTime := now""" + $newline +
    "    - context: ""timeDiff := timestamp(context.currentTime) - timestamp(context.lastLogin)""" + $newline +
    "    - if:" + $newline +
    "        match: context.timeDiff <= duration('30m')" + $newline +
    "        block:" + $newline +
    "        - return: true" + $newline +
    "    - return: false"

$targetRow = 69

$ws.Range("A$targetRow").Value = $aText
$ws.Range("B$targetRow").Value = $bText

# Same wrap-text style ("s=1") used by every other row in the table.
$rowRange = $ws.Range("A" + $targetRow + ":B" + $targetRow)
$rowRange.WrapText = $true

# Row grows tall enough to show all nine wrapped lines of the rule body.
$ws.Rows.Item($targetRow).RowHeight = 129.6
